$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G (K) values per row, as recomputed from play-by-play source (K instead of Strike#).
$gValues = @{
    2 = 0
    3 = 2
    4 = 0
    6 = 0
    7 = 5
    8 = 2
    9 = 2
    10 = 0
    11 = 1
    12 = 2
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 2
    20 = 4
    21 = 2
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 2
    27 = 2
    28 = 2
    29 = 1
    30 = 1
    31 = 2
    32 = 1
    33 = 0
    34 = 1
    36 = 2
    37 = 0
    38 = 0
    39 = 1
    40 = 0
    42 = 2
    43 = 0
    44 = 3
    46 = 3
    47 = 1
    48 = 0
    49 = 2
    50 = 0
    51 = 0
    52 = 1
    53 = 2
    54 = 1
    55 = 1
    56 = 2
    58 = 1
    59 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
